# Apply daily update to "Croatia HNL" odds dataset:
#  - rows 131 and 132 (matches that have since been played) receive their
#    full-time result (FTHG/FTAG/FTR) and closing odds/profit-loss columns
#  - three new fixture rows (133-135) are appended, one with a result already
#    available (row 133) and two still pending (rows 134-135)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 131 (id 129, NK Osijek v NK Varazdin): add result + refresh closing odds ---
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 1
$ws.Range("J131").Value = "A"
$ws.Range("N131").Value = 1.666
$ws.Range("O131").Value = 3.6
$ws.Range("P131").Value = 4.75
$ws.Range("Q131").Value = -0.75
$ws.Range("R131").Value = 1.925
$ws.Range("S131").Value = 1.925
$ws.Range("U131").Value = 2.025
$ws.Range("V131").Value = 1.825
$ws.Range("W131").Value = -1
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = 3.75
$ws.Range("Z131").Value = -1
$ws.Range("AA131").Value = 0.925
$ws.Range("AB131").Value = -1
$ws.Range("AC131").Value = 0.825

# --- Row 132 (id 130, now the Slaven Belupo v HNK Rijeka match): full refresh ---
$ws.Range("B132").Value = 6788934
$ws.Range("E132").Value = 45367.45833333334
$ws.Range("F132").Value = "Slaven Belupo"
$ws.Range("G132").Value = "HNK Rijeka"
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 1
$ws.Range("J132").Value = "A"
$ws.Range("K132").Value = 5.5
$ws.Range("L132").Value = 3.8
$ws.Range("M132").Value = 1.55
$ws.Range("N132").Value = 6
$ws.Range("O132").Value = 3.8
$ws.Range("P132").Value = 1.5
$ws.Range("Q132").Value = 1
$ws.Range("R132").Value = 2
$ws.Range("S132").Value = 1.85
$ws.Range("T132").Value = 2.5
$ws.Range("W132").Value = -1
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = 0.5
$ws.Range("AA132").Value = 0
$ws.Range("AB132").Value = -1
$ws.Range("AC132").Value = 0.875

# --- New row 133 (id 131) ---
$ws.Range("A130").Copy()
$ws.Range("A133").PasteSpecial(-4122)
$ws.Range("E130").Copy()
$ws.Range("E133").PasteSpecial(-4122)
$ws.Range("A133").Value = 131
$ws.Range("B133").Value = 6787895
$ws.Range("C133").Value = "Croatia HNL"
$ws.Range("D133").Value = "Croatia HNL"
$ws.Range("E133").Value = 45367.54861111111
$ws.Range("F133").Value = "Hajduk Split"
$ws.Range("G133").Value = "NK Lokomotiva Zagreb"
$ws.Range("H133").Value = 1
$ws.Range("I133").Value = 2
$ws.Range("J133").Value = "A"
$ws.Range("K133").Value = 1.363
$ws.Range("L133").Value = 4.2
$ws.Range("M133").Value = 8.5
$ws.Range("N133").Value = 1.363
$ws.Range("O133").Value = 4.2
$ws.Range("P133").Value = 8.5
$ws.Range("Q133").Value = -1.25
$ws.Range("R133").Value = 1.875
$ws.Range("S133").Value = 1.975
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.925
$ws.Range("V133").Value = 1.925
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = -1
$ws.Range("Y133").Value = 7.5
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.9750000000000001
$ws.Range("AB133").Value = 0.925
$ws.Range("AC133").Value = -1

# --- New row 134 (id 132) ---
$ws.Range("A130").Copy()
$ws.Range("A134").PasteSpecial(-4122)
$ws.Range("E130").Copy()
$ws.Range("E134").PasteSpecial(-4122)
$ws.Range("A134").Value = 132
$ws.Range("B134").Value = 6788935
$ws.Range("C134").Value = "Croatia HNL"
$ws.Range("D134").Value = "Croatia HNL"
$ws.Range("E134").Value = 45368.45833333334
$ws.Range("F134").Value = "Istra 1961"
$ws.Range("G134").Value = "HNK Gorica"
$ws.Range("K134").Value = 2
$ws.Range("L134").Value = 3.2
$ws.Range("M134").Value = 3.8
$ws.Range("N134").Value = 1.909
$ws.Range("O134").Value = 3.2
$ws.Range("P134").Value = 4.333
$ws.Range("Q134").Value = -0.5
$ws.Range("R134").Value = 1.95
$ws.Range("S134").Value = 1.9
$ws.Range("T134").Value = 2.25
$ws.Range("U134").Value = 2.05
$ws.Range("V134").Value = 1.8
$ws.Range("W134").Value = 0
$ws.Range("X134").Value = 0
$ws.Range("Y134").Value = 0
$ws.Range("Z134").Value = 0
$ws.Range("AA134").Value = 0

# --- New row 135 (id 133) ---
$ws.Range("A130").Copy()
$ws.Range("A135").PasteSpecial(-4122)
$ws.Range("E130").Copy()
$ws.Range("E135").PasteSpecial(-4122)
$ws.Range("A135").Value = 133
$ws.Range("B135").Value = 6769305
$ws.Range("C135").Value = "Croatia HNL"
$ws.Range("D135").Value = "Croatia HNL"
$ws.Range("E135").Value = 45368.54861111111
$ws.Range("F135").Value = "NK Rudes"
$ws.Range("G135").Value = "Dinamo Zagreb"
$ws.Range("K135").Value = 15
$ws.Range("L135").Value = 6
$ws.Range("M135").Value = 1.166
$ws.Range("N135").Value = 15
$ws.Range("O135").Value = 7
$ws.Range("P135").Value = 1.142
$ws.Range("Q135").Value = 2
$ws.Range("R135").Value = 1.975
$ws.Range("S135").Value = 1.875
$ws.Range("T135").Value = 2.75
$ws.Range("U135").Value = 1.8
$ws.Range("V135").Value = 2.05
$ws.Range("W135").Value = 0
$ws.Range("X135").Value = 0
$ws.Range("Y135").Value = 0
$ws.Range("Z135").Value = 0
$ws.Range("AA135").Value = 0

